$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row height tweaks ---
$ws.Rows.Item(2).RowHeight = 14.15
$ws.Rows.Item(3).RowHeight = 14.15
$ws.Rows.Item(4).RowHeight = 14.15

# --- Plain text / number value updates shared across rows 2-5 ---
$ws.Range("E2:E5").Value = "Maandelijkse Facturatie a €1000,= geschatte éénmalige investeringskosten n.v.t.Y = aantal jaren looptijd overeenkomstX = geschatte jaarlijkse kosten"
$ws.Range("G2:G4").Value = "Eelco Aartsen"
$ws.Range("I2:I5").Value = "075-6163455"
$ws.Range("R2:R5").Value = 2300
$ws.Range("T2:T5").Value = "Pietje Puk"
$ws.Range("V2:V5").Value = "06-1231231"
$ws.Range("Z2:Z5").Value = "Sneeuwbezems 2021"

# --- Hyperlinks: drop the old H5 one, re-add everything so the
#     Contracteigenaar (H) and Contactpersoon (U) e-mail columns all
#     point at the new addresses with matching blue hyperlink styling ---
$ws.Range("H5").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:eelco@aesset.nl", "", "", "eelco@aesset.nl") | Out-Null
$ws.Hyperlinks.Add($ws.Range("U2"), "mailto:p.puk@npo.nl", "", "", "p.puk@npo.nl") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H3"), "mailto:eelco@aesset.nl", "", "", "eelco@aesset.nl") | Out-Null
$ws.Hyperlinks.Add($ws.Range("U3"), "mailto:p.puk@npo.nl", "", "", "p.puk@npo.nl") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H4"), "mailto:eelco@aesset.nl", "", "", "eelco@aesset.nl") | Out-Null
$ws.Hyperlinks.Add($ws.Range("U4"), "mailto:p.puk@npo.nl", "", "", "p.puk@npo.nl") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H5"), "mailto:eelco@aesset.nl", "", "", "eelco@aesset.nl") | Out-Null
$ws.Hyperlinks.Add($ws.Range("U5"), "mailto:p.puk@npo.nl", "", "", "p.puk@npo.nl") | Out-Null

# Hyperlinks.Add stamps Excel's default "Hyperlink" style (blue + underline);
# re-apply the workbook's existing blue, non-underlined look (same as the
# original H5 contracteigenaar cell) across every touched hyperlink cell so
# no new font/style gets minted.
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H2,H3,H4,H5,U2,U3,U4,U5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
